$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.808.50"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "3.898.43"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "475.81"
$ws.Range("E5").Value = "  +5.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.10"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +8.42%  "

$ws.Range("E11").Value = "  +11.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.37"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").Value = "4.499.47"
$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.25"
$ws.Range("E14").Value = "  -0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.59"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("D16").Value = "3.940.14"
$ws.Range("E16").Value = "  +2.77%  "

$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.75"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("E19").Value = "  -3.66%  "

$ws.Range("D20").Value = "67.710.02"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.44"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.33"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.31"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.96"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.10"
$ws.Range("E26").Value = "  +2.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  -0.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  +3.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "723.72"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.22"
$ws.Range("E30").Value = "  -4.08%  "

$ws.Range("E31").Value = "  -4.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "41.95"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0887"
$ws.Range("E34").Value = "  +30.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.74"
$ws.Range("E35").Value = "  +2.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").Value = "  -4.39%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.36"
$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0467"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  +6.59%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "30.20"
$ws.Range("E41").Value = "  +21.11%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.98"
$ws.Range("E42").Value = "  +2.85%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("E43").Value = "  +11.28%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.340"
$ws.Range("E44").Value = "  -2.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.140"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.18"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.47"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.86"
$ws.Range("E51").Value = "  -0.02%  "
